$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay plain text (matches source data which
# includes values like "27.643.45" / "1.029" that Excel would otherwise
# coerce into numbers or dates).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.643.45"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "1.853.97"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("D4").Value = "1.029"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "321.30"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").Value = "1.025"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "0.4377"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  +2.05%  "
$ws.Range("D9").Value = "0.07415"
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").Value = "0.8816"
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("D11").Value = "21.58"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").Value = "1.860.65"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "5.518"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").Value = "6.725"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").Value = "0.07147"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "85.21"
$ws.Range("E16").Value = "  +3.27%  "
$ws.Range("D17").Value = "1.029"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "0.000009095"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").Value = "1.025"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "15.47"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D21").Value = "27.661.38"
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").Value = "5.296"
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").Value = "11.26"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "2.090.49"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").Value = "2.011"
$ws.Range("E25").Value = "  +5.33%  "
$ws.Range("D26").Value = "157.30"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").Value = "18.76"
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("D28").Value = "5.348"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").Value = "1.986"
$ws.Range("E29").Value = "  +3.45%  "
$ws.Range("D30").Value = "117.74"
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("D31").Value = "0.08992"
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").Value = "0.7772"
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("D33").Value = "1.214"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").Value = "2.990"
$ws.Range("E34").Value = "  +4.42%  "
$ws.Range("D35").Value = "4.562"
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("D36").Value = "1.027"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").Value = "1.142"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("D38").Value = "0.01973"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "0.05275"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").Value = "2.863"
$ws.Range("E40").Value = "  +2.99%  "
$ws.Range("D41").Value = "0.5183"
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("D42").Value = "0.1681"
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("D43").Value = "6.852"
$ws.Range("E43").Value = "  +4.67%  "
$ws.Range("D44").Value = "8.871"
$ws.Range("E44").Value = "  +4.56%  "
$ws.Range("D45").Value = "110.01"
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("D46").Value = "10.73"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("D47").Value = "0.06612"
$ws.Range("E47").Value = "  +4.69%  "
$ws.Range("D48").Value = "1.027"
$ws.Range("E48").Value = "  +0.18%  "
$ws.Range("D49").Value = "1.705"
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("D50").Value = "0.4710"
$ws.Range("E50").Value = "  +1.66%  "
$ws.Range("E51").Value = "  -0.33%  "
